# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H, I, J, K, L, M, N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# with latest market-board snapshot values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1635.64
$ws.Range("I19").Value = 902.25
$ws.Range("J19").Value = 2939.4443
$ws.Range("K19").Value = 902.25
$ws.Range("L19").Value = 2939.4443
$ws.Range("M19").Value = -727.25
$ws.Range("N19").Value = -3289.4443
$ws.Range("H21").Value = 10411.765
$ws.Range("H23").Value = 10411.765
$ws.Range("H29").Value = 4120
$ws.Range("I29").Value = 240
$ws.Range("J29").Value = 8000
$ws.Range("K29").Value = 720
$ws.Range("L29").Value = 24000
$ws.Range("M29").Value = -439
$ws.Range("N29").Value = -24562
$ws.Range("H38").Value = 5960.3335
$ws.Range("I38").Value = 127.375
$ws.Range("J38").Value = 9549.846
$ws.Range("K38").Value = 382.125
$ws.Range("L38").Value = 28649.538
$ws.Range("M38").Value = -10.125
$ws.Range("N38").Value = -29393.538
$ws.Range("H64").Value = 7318.1816
$ws.Range("J64").Value = 8666.666999999999
$ws.Range("L64").Value = 8666.666999999999
$ws.Range("N64").Value = -9162.666999999999
$ws.Range("H67").Value = 7318.1816
$ws.Range("J67").Value = 8666.666999999999
$ws.Range("L67").Value = 8666.666999999999
$ws.Range("N67").Value = -10382.667
$ws.Range("H80").Value = 1029.5161
$ws.Range("I80").Value = 1021.5
$ws.Range("J80").Value = 1033.3334
$ws.Range("K80").Value = 3064.5
$ws.Range("L80").Value = 3100.0002
$ws.Range("M80").Value = -2066.5
$ws.Range("N80").Value = -5096.0002
$ws.Range("H83").Value = 1029.5161
$ws.Range("I83").Value = 1021.5
$ws.Range("J83").Value = 1033.3334
$ws.Range("K83").Value = 9193.5
$ws.Range("L83").Value = 9300.000599999999
$ws.Range("M83").Value = -4201.5
$ws.Range("N83").Value = -19284.0006
$ws.Range("H125").Value = 2036
$ws.Range("J125").Value = 2036
$ws.Range("L125").Value = 18324
$ws.Range("N125").Value = -23244
$ws.Range("H135").Value = 2855.9285
$ws.Range("I135").Value = 3613.4443
$ws.Range("J135").Value = 1492.4
$ws.Range("K135").Value = 32520.9987
$ws.Range("L135").Value = 13431.6
$ws.Range("M135").Value = -29985.9987
$ws.Range("N135").Value = -18501.6
$ws.Range("H139").Value = 99994.2
$ws.Range("J139").Value = 99994.2
$ws.Range("L139").Value = 99994.2
$ws.Range("N139").Value = -110274.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 49709.848
$ws.Range("I102").Value = 12185.75
$ws.Range("K102").Value = 12185.75
$ws.Range("M102").Value = -10563.75
$ws.Range("H132").Value = 3872.0566
$ws.Range("J132").Value = 4260
$ws.Range("L132").Value = 12780
$ws.Range("N132").Value = -17840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 37500
$ws.Range("J117").Value = 37500
$ws.Range("L117").Value = 37500
$ws.Range("N117").Value = -46678
$ws.Range("H134").Value = 1552.093
$ws.Range("I134").Value = 1459.2439
$ws.Range("J134").Value = 3455.5
$ws.Range("K134").Value = 4377.7317
$ws.Range("L134").Value = 10366.5
$ws.Range("M134").Value = -1842.7317
$ws.Range("N134").Value = -15436.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2039.2759
$ws.Range("J31").Value = 2300.7693
$ws.Range("L31").Value = 2300.7693
$ws.Range("N31").Value = -2890.7693
$ws.Range("H34").Value = 2039.2759
$ws.Range("J34").Value = 2300.7693
$ws.Range("L34").Value = 2300.7693
$ws.Range("N34").Value = -2704.7693
$ws.Range("H132").Value = 1825.1621
$ws.Range("I132").Value = 1700.8857
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5102.6571
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2572.6571
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 150.26666
$ws.Range("I12").Value = 86.5
$ws.Range("J12").Value = 223.14285
$ws.Range("K12").Value = 259.5
$ws.Range("L12").Value = 669.4285500000001
$ws.Range("M12").Value = -86.5
$ws.Range("N12").Value = -1015.42855
$ws.Range("H38").Value = 116.875
$ws.Range("I38").Value = 116
$ws.Range("K38").Value = 348
$ws.Range("M38").Value = -1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 45743.75
$ws.Range("J136").Value = 45743.75
$ws.Range("L136").Value = 137231.25
$ws.Range("N136").Value = -142331.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1384.2858
$ws.Range("I31").Value = 1281.6666
$ws.Range("K31").Value = 1281.6666
$ws.Range("M31").Value = -1033.6666
$ws.Range("H82").Value = 1600.3077
$ws.Range("I82").Value = 1535.625
$ws.Range("J82").Value = 1703.8
$ws.Range("K82").Value = 1535.625
$ws.Range("L82").Value = 1703.8
$ws.Range("M82").Value = -1174.625
$ws.Range("N82").Value = -2425.8
$ws.Range("H85").Value = 1600.3077
$ws.Range("I85").Value = 1535.625
$ws.Range("J85").Value = 1703.8
$ws.Range("K85").Value = 1535.625
$ws.Range("L85").Value = 1703.8
$ws.Range("M85").Value = -287.625
$ws.Range("N85").Value = -4199.8
$ws.Range("H100").Value = 1670.2084
$ws.Range("I100").Value = 871.9524
$ws.Range("K100").Value = 871.9524
$ws.Range("M100").Value = -330.9524
$ws.Range("H130").Value = 49300
$ws.Range("J130").Value = 49300
$ws.Range("L130").Value = 49300
$ws.Range("N130").Value = -59340
$ws.Range("H132").Value = 3938.7083
$ws.Range("I132").Value = 4018.1667
$ws.Range("J132").Value = 3700.3333
$ws.Range("K132").Value = 12054.5001
$ws.Range("L132").Value = 11100.9999
$ws.Range("M132").Value = -9524.500100000001
$ws.Range("N132").Value = -16160.9999
$ws.Range("H133").Value = 65633.336
$ws.Range("J133").Value = 65633.336
$ws.Range("L133").Value = 65633.336
$ws.Range("N133").Value = -70693.336
$ws.Range("H136").Value = 3530.1333
$ws.Range("I136").Value = 3231.5454
$ws.Range("J136").Value = 4351.25
$ws.Range("K136").Value = 9694.636200000001
$ws.Range("L136").Value = 13053.75
$ws.Range("M136").Value = -7144.636200000001
$ws.Range("N136").Value = -18153.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 115162.5
$ws.Range("J86").Value = 115162.5
$ws.Range("L86").Value = 115162.5
$ws.Range("N86").Value = -117408.5
$ws.Range("H89").Value = 115162.5
$ws.Range("J89").Value = 115162.5
$ws.Range("L89").Value = 575812.5
$ws.Range("N89").Value = -587044.5

